$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 141, shifting existing rows 141:189 down to 142:190
$ws.Rows("141:141").Insert()

# Populate the newly inserted row 141 with the new weekly record
$ws.Range("A141").Value = 3
$ws.Range("B141").Value = "Femacal de La Calera"
$ws.Range("C141").Value = "Coquimbo"
$ws.Range("D141").Value = 44704
$ws.Range("E141").Value = 5
$ws.Range("F141").Value = 100112030
$ws.Range("G141").Value = "Poroto granado"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 76
$ws.Range("K141").Value = 23000
$ws.Range("L141").Value = 24000
$ws.Range("M141").Value = 23500
$ws.Range("N141").Value = "$/saco 25 kilos"
$ws.Range("O141").Value = "Provincia de Limarí"
$ws.Range("P141").Value = 940
$ws.Range("Q141").Value = 25
$ws.Range("R141").Value = "Hortaliza"

# Preserve the date-column number format for the new row (matches style of column D elsewhere)
$ws.Range("D141").NumberFormat = "YYYY-MM-DD HH:MM:SS"
